$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.615317344665527
$ws.Range("B1").Value = 6.209388732910156
$ws.Range("C1").Value = 5.602262020111084
$ws.Range("D1").Value = 6.432198047637939
$ws.Range("E1").Value = 3.850381851196289
